$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (Text) number format on cells whose new values would
# otherwise be auto-converted to numbers by Excel (which would lose
# literal formatting such as trailing zeros, e.g. "0.530" or "1.00").
# NumberFormat must be applied per-cell because it only affects the
# first area of a multi-area (union) range.
$textCells = @("D5", "D6", "D8", "D10", "D11", "D12", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.119.96"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "3.452.03"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "584.62"
$ws.Range("E5").Value = "  -2.65%  "
$ws.Range("D6").Value = "176.02"
$ws.Range("E6").Value = "  -2.91%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.608"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").Value = "3.452.18"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").Value = "0.135"
$ws.Range("E10").Value = "  -3.87%  "
$ws.Range("D11").Value = "6.93"
$ws.Range("E11").Value = "  -3.11%  "
$ws.Range("D12").Value = "0.424"
$ws.Range("E12").Value = "  -3.18%  "
$ws.Range("D13").Value = "4.064.35"
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").Value = "31.51"
$ws.Range("E14").Value = "  -3.78%  "
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").Value = "67.141.55"
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("E17").Value = "  -3.63%  "
$ws.Range("D18").Value = "3.440.23"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").Value = "6.07"
$ws.Range("E19").Value = "  -4.62%  "
$ws.Range("D20").Value = "13.87"
$ws.Range("E20").Value = "  -4.14%  "
$ws.Range("D21").Value = "379.00"
$ws.Range("E21").Value = "  -5.23%  "
$ws.Range("D22").Value = "7.81"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("D23").Value = "5.77"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.16%  "
$ws.Range("D25").Value = "71.57"
$ws.Range("E25").Value = "  -2.86%  "
$ws.Range("D26").Value = "0.530"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").Value = "9.87"
$ws.Range("E28").Value = "  -6.14%  "
$ws.Range("D29").Value = "0.173"
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").Value = "5.95"
$ws.Range("E31").Value = "  -4.82%  "
$ws.Range("D32").Value = "2.02"
$ws.Range("E32").Value = "  -3.46%  "
$ws.Range("D33").Value = "23.94"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("E34").Value = "  -6.48%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("D36").Value = "7.19"
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("D37").Value = "1.54"
$ws.Range("E37").Value = "  -4.98%  "
$ws.Range("D38").Value = "160.12"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").Value = "0.877"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").Value = "26.89"
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("D41").Value = "1.81"
$ws.Range("E41").Value = "  -5.85%  "
$ws.Range("D42").Value = "2.63"
$ws.Range("E42").Value = "  -4.78%  "
$ws.Range("E43").Value = "  -5.11%  "
$ws.Range("D44").Value = "4.47"
$ws.Range("E44").Value = "  -4.63%  "
$ws.Range("D45").Value = "2.695.14"
$ws.Range("E45").Value = "  -6.76%  "
$ws.Range("D46").Value = "0.0696"
$ws.Range("E46").Value = "  -5.38%  "
$ws.Range("D47").Value = "25.39"
$ws.Range("E47").Value = "  -5.63%  "
$ws.Range("D48").Value = "41.08"
$ws.Range("E48").Value = "  -3.08%  "
$ws.Range("D49").Value = "0.0294"
$ws.Range("E49").Value = "  -3.21%  "
$ws.Range("D50").Value = "321.19"
$ws.Range("E50").Value = "  -7.06%  "
$ws.Range("D51").Value = "1.02"
$ws.Range("E51").Value = "  -4.20%  "

Write-Output "Applied cryptos update."
